$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.699.26"
Set-TextValue "E2" "  +1.51%  "
Set-TextValue "D3" "1.637.02"
Set-TextValue "E3" "  +0.95%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "213.29"
Set-TextValue "E5" "  +0.48%  "
Set-TextValue "D6" "0.502"
Set-TextValue "E6" "  +2.93%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "E8" "  +2.65%  "
Set-TextValue "E9" "  +1.33%  "
Set-TextValue "D10" "19.24"
Set-TextValue "E10" "  +1.74%  "
Set-TextValue "D11" "0.0842"
Set-TextValue "E11" "  +3.43%  "
Set-TextValue "D12" "1.866.10"
Set-TextValue "E12" "  +0.97%  "
Set-TextValue "D13" "1.648.55"
Set-TextValue "E13" "  +1.61%  "
Set-TextValue "E14" "  +2.27%  "
Set-TextValue "E15" "  +1.77%  "
Set-TextValue "D16" "26.693.92"
Set-TextValue "E16" "  +1.42%  "
Set-TextValue "D17" "63.31"
Set-TextValue "D18" "0.0₃0746"
Set-TextValue "E18" "  +2.33%  "
Set-TextValue "D19" "217.95"
Set-TextValue "E19" "  +7.51%  "
Set-TextValue "E20" "  +0.04%  "
Set-TextValue "E21" "  +0.46%  "
Set-TextValue "E22" "  +1.20%  "
Set-TextValue "D23" "6.22"
Set-TextValue "E23" "  +2.72%  "
Set-TextValue "E24" "  +0.36%  "
Set-TextValue "D25" "148.39"
Set-TextValue "E25" "  +3.04%  "
Set-TextValue "E26" "  +0.02%  "
Set-TextValue "E27" "  +0.26%  "
Set-TextValue "D28" "7.01"
Set-TextValue "E28" "  +6.74%  "
Set-TextValue "D29" "15.46"
Set-TextValue "E29" "  +1.64%  "
Set-TextValue "D30" "0.0510"
Set-TextValue "E30" "  -3.48%  "
Set-TextValue "E31" "  -0.21%  "
Set-TextValue "D32" "3.33"
Set-TextValue "E32" "  +4.56%  "
Set-TextValue "D33" "2.96"
Set-TextValue "E33" "  +0.37%  "
Set-TextValue "E34" "  +0.94%  "
Set-TextValue "E35" "  -0.89%  "
Set-TextValue "D36" "1.205.00"
Set-TextValue "E36" "  +1.95%  "
Set-TextValue "E37" "  +6.05%  "
Set-TextValue "D38" "0.811"
Set-TextValue "E38" "  +0.25%  "
Set-TextValue "E39" "  +0.03%  "
Set-TextValue "E40" "  +1.91%  "
Set-TextValue "D41" "2.30"
Set-TextValue "E41" "  -0.90%  "
Set-TextValue "D42" "5.43"
Set-TextValue "E42" "  +1.50%  "
Set-TextValue "E43" "  +0.67%  "
Set-TextValue "D44" "1.773.91"
Set-TextValue "E44" "  +0.83%  "
Set-TextValue "D45" "92.56"
Set-TextValue "E45" "  -0.97%  "
Set-TextValue "E46" "  +2.18%  "
Set-TextValue "D47" "54.93"
Set-TextValue "E47" "  +1.86%  "
Set-TextValue "D48" "0.0513"
Set-TextValue "E48" "  +0.82%  "
Set-TextValue "D49" "7.67"
Set-TextValue "E49" "  +5.69%  "
Set-TextValue "D50" "0.411"
Set-TextValue "E50" "  +0.53%  "
Set-TextValue "E51" "  +0.03%  "
